$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets("ALC")
# Row 132
$ws.Range("H132").Value = 22553594
$ws.Range("I132").Value = 28302238
$ws.Range("J132").Value = 1222.7693
$ws.Range("K132").Value = 84906714
$ws.Range("L132").Value = 3668.3079
$ws.Range("M132").Value = -84904184
$ws.Range("N132").Value = -8728.3079

# Row 137
$ws.Range("H137").Value = 1314.8572
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 1314.8572
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 3944.5716
$ws.Range("N137").Value = -9044.571599999999
$ws.Range("M137").ClearContents()

# Row 138
$ws.Range("H138").Value = 2383.9456
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 2383.9456
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 7151.8368
$ws.Range("N138").Value = -17431.8368
$ws.Range("M138").ClearContents()


# --- Sheet: ARM ---
$ws = $wb.Worksheets("ARM")
# Row 32
$ws.Range("H32").Value = 5337.879
$ws.Range("I32").Value = 4311.4873
$ws.Range("K32").Value = 4311.4873
$ws.Range("M32").Value = -4024.4873

# Row 45
$ws.Range("H45").Value = 873885.6
$ws.Range("I45").Value = 1209651.6
$ws.Range("J45").Value = 894
$ws.Range("K45").Value = 1209651.6
$ws.Range("L45").Value = 894
$ws.Range("M45").Value = -1209274.6
$ws.Range("N45").Value = -1648

# Row 61
$ws.Range("H61").Value = 3612.8708
$ws.Range("I61").Value = 4156.36
$ws.Range("J61").Value = 1348.3334
$ws.Range("K61").Value = 4156.36
$ws.Range("L61").Value = 1348.3334
$ws.Range("M61").Value = -3944.36
$ws.Range("N61").Value = -1772.3334

# Row 132
$ws.Range("H132").Value = 2366.25
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 2366.25
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 7098.75
$ws.Range("N132").Value = -12158.75
$ws.Range("M132").ClearContents()

# Row 136
$ws.Range("H136").Value = 3612.8708
$ws.Range("I136").Value = 4156.36
$ws.Range("J136").Value = 1348.3334
$ws.Range("K136").Value = 12469.08
$ws.Range("L136").Value = 4045.0002
$ws.Range("M136").Value = -9919.079999999998
$ws.Range("N136").Value = -9145.0002


# --- Sheet: BSM ---
$ws = $wb.Worksheets("BSM")
# Row 134
$ws.Range("H134").Value = 10769742
$ws.Range("I134").Value = 13910421
$ws.Range("J134").Value = 1702
$ws.Range("K134").Value = 41731263
$ws.Range("L134").Value = 5106
$ws.Range("M134").Value = -41728728
$ws.Range("N134").Value = -10176


# --- Sheet: CRP ---
$ws = $wb.Worksheets("CRP")
# Row 31
$ws.Range("H31").Value = 21746968
$ws.Range("I31").Value = 26316862
$ws.Range("J31").Value = 39967.75
$ws.Range("K31").Value = 26316862
$ws.Range("L31").Value = 39967.75
$ws.Range("M31").Value = -26316567
$ws.Range("N31").Value = -40557.75

# Row 34
$ws.Range("H34").Value = 21746968
$ws.Range("I34").Value = 26316862
$ws.Range("J34").Value = 39967.75
$ws.Range("K34").Value = 26316862
$ws.Range("L34").Value = 39967.75
$ws.Range("M34").Value = -26316660
$ws.Range("N34").Value = -40371.75

# Row 58
$ws.Range("H58").Value = 2999911.8
$ws.Range("I58").Value = 3197596.8
$ws.Range("J58").Value = 34638
$ws.Range("K58").Value = 3197596.8
$ws.Range("L58").Value = 34638
$ws.Range("M58").Value = -3197393.8
$ws.Range("N58").Value = -35044

# Row 132
$ws.Range("H132").Value = 5955851
$ws.Range("I132").Value = 8334916
$ws.Range("J132").Value = 8188.875
$ws.Range("K132").Value = 25004748
$ws.Range("L132").Value = 24566.625
$ws.Range("M132").Value = -25002218
$ws.Range("N132").Value = -29626.625

# Row 134
$ws.Range("H134").Value = 31251348
$ws.Range("I134").Value = 43104550
$ws.Range("J134").Value = 4809586.5
$ws.Range("K134").Value = 129313650
$ws.Range("L134").Value = 14428759.5
$ws.Range("M134").Value = -129311115
$ws.Range("N134").Value = -14433829.5

# Row 136
$ws.Range("H136").Value = 2999911.8
$ws.Range("I136").Value = 3197596.8
$ws.Range("J136").Value = 34638
$ws.Range("K136").Value = 9592790.399999999
$ws.Range("L136").Value = 103914
$ws.Range("M136").Value = -9590240.399999999
$ws.Range("N136").Value = -109014


# --- Sheet: CUL ---
$ws = $wb.Worksheets("CUL")
# Row 5
$ws.Range("H5").Value = 657.5333000000001
$ws.Range("I5").Value = 573.5
$ws.Range("J5").Value = 688.0909
$ws.Range("K5").Value = 1720.5
$ws.Range("L5").Value = 2064.2727
$ws.Range("M5").Value = -1608.5
$ws.Range("N5").Value = -2288.2727

# Row 135
$ws.Range("H135").Value = 657.5333000000001
$ws.Range("I135").Value = 573.5
$ws.Range("J135").Value = 688.0909
$ws.Range("K135").Value = 5161.5
$ws.Range("L135").Value = 6192.8181
$ws.Range("M135").Value = -2626.5
$ws.Range("N135").Value = -11262.8181


# --- Sheet: GSM ---
$ws = $wb.Worksheets("GSM")
# Row 113
$ws.Range("H113").Value = 1185.5454
$ws.Range("I113").Value = 905.8570999999999
$ws.Range("J113").Value = 1675
$ws.Range("K113").Value = 905.8570999999999
$ws.Range("L113").Value = 1675
$ws.Range("M113").Value = 1264.1429
$ws.Range("N113").Value = -6015

# Row 132
$ws.Range("H132").Value = 39218828
$ws.Range("I132").Value = 64517556
$ws.Range("J132").Value = 5798.6
$ws.Range("K132").Value = 193552668
$ws.Range("L132").Value = 17395.8
$ws.Range("M132").Value = -193550138
$ws.Range("N132").Value = -22455.8


# --- Sheet: LTW ---
$ws = $wb.Worksheets("LTW")
# Row 43
$ws.Range("H43").Value = 5698.5
$ws.Range("J43").Value = 5698.5
$ws.Range("L43").Value = 5698.5
$ws.Range("N43").Value = -6084.5

# Row 61
$ws.Range("H61").Value = 1477.3572
$ws.Range("I61").Value = 1061.1818
$ws.Range("J61").Value = 3003.3333
$ws.Range("K61").Value = 1061.1818
$ws.Range("L61").Value = 3003.3333
$ws.Range("M61").Value = -859.1818000000001
$ws.Range("N61").Value = -3407.3333

# Row 97
$ws.Range("H97").Value = 0
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()

# Row 103
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Row 113
$ws.Range("H113").Value = 1477.3572
$ws.Range("I113").Value = 1061.1818
$ws.Range("J113").Value = 3003.3333
$ws.Range("K113").Value = 1061.1818
$ws.Range("L113").Value = 3003.3333
$ws.Range("M113").Value = 1108.8182
$ws.Range("N113").Value = -7343.3333

# Row 132
$ws.Range("H132").Value = 5130075.5
$ws.Range("I132").Value = 7693678.5
$ws.Range("J132").Value = 2869.923
$ws.Range("K132").Value = 23081035.5
$ws.Range("L132").Value = 8609.769
$ws.Range("M132").Value = -23078505.5
$ws.Range("N132").Value = -13669.769

# Row 136
$ws.Range("H136").Value = 3598.392
$ws.Range("I136").Value = 3815.6
$ws.Range("J136").Value = 1969.3334
$ws.Range("K136").Value = 11446.8
$ws.Range("L136").Value = 5908.0002
$ws.Range("M136").Value = -8896.799999999999
$ws.Range("N136").Value = -11008.0002


# --- Sheet: WVR ---
$ws = $wb.Worksheets("WVR")
# Row 113
$ws.Range("H113").Value = 858.6087
$ws.Range("I113").Value = 559
$ws.Range("J113").Value = 1185.4546
$ws.Range("K113").Value = 1677
$ws.Range("L113").Value = 3556.3638
$ws.Range("M113").Value = 493
$ws.Range("N113").Value = -7896.3638

# Row 132
$ws.Range("H132").Value = 7678766
$ws.Range("I132").Value = 4652206
$ws.Range("J132").Value = 10281607
$ws.Range("K132").Value = 13956618
$ws.Range("L132").Value = 30844821
$ws.Range("M132").Value = -13954088
$ws.Range("N132").Value = -30849881

